# ---------------------------------------------------------------------------
# Edits "Storm 200" listing-scrape document to match the target revision:
#   1. A new "EVOX ELETRONICOS" listing block is inserted at the very start
#      of the document (duplicate data of the block that used to sit right
#      after "KGMICOMERCIAL"; that old occurrence is removed later).
#   2. The "RADICALSOM.COM.BR" listing is split in two:
#        - a brand new block (different MLB id, price, "Premium" tipo)
#          inserted just before the original one,
#        - the original block's MLB id changes, and its Loja/Lugar move to
#          a new branch ("RADICAL_SOM_FILIAL" / "Joinville, Santa
#          Catarina.").
#   3. The old, now-duplicated "EVOX ELETRONICOS" block (the one that used
#      to directly follow "KGMICOMERCIAL") is deleted.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# Helper: 1-based paragraph index of the paragraph that contains the first
# occurrence of $searchText in the whole document.
function Get-ParagraphIndexContaining([string]$searchText) {
    $hit = $d.Content
    $hit.Find.ClearFormatting()
    $found = $hit.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $searchText"
    }
    $prefix = $d.Range(0, $hit.End)
    return $prefix.Paragraphs.Count
}

# Helper: build the CRLF-joined text (using Word's paragraph-mark escape,
# "`r") for one full 9-line listing block plus its trailing blank separator
# paragraph, ready to be handed to Range.InsertBefore / InsertAfter.
function New-BlockText([string]$modelo, [string]$url, [string]$nome, [string]$preco,
                        [string]$precoPrevisto, [string]$loja, [string]$tipo, [string]$lugar) {
    $lines = @(
        "Modelo: $modelo",
        "URL: $url",
        "Nome: $nome",
        "Preço: $preco",
        "Preço Previsto: $precoPrevisto",
        "Loja: $loja",
        "Tipo: $tipo",
        "Lugar: $lugar",
        "--------------------------------------------------------------------"
    )
    return ($lines -join "`r") + "`r`r"
}

# ---------------------------------------------------------------------------
# Step 1: rename the MLB id of the RADICALSOM.COM.BR listing
#         (3452387362 -> 4156504566). The id is unique in the doc at this
#         point, so a plain Find/Replace is unambiguous.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "MLB-3452387362-carregador-jfa-storm-200a-144v-mais-completa-smart-cca-220v-_JM",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "MLB-4156504566-carregador-jfa-storm-200a-144v-mais-completa-smart-cca-220v-_JM", 2) | Out-Null

# ---------------------------------------------------------------------------
# Step 2: that same (original) listing moves to a new branch.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Loja: RADICALSOM.COM.BR", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Loja: RADICAL_SOM_FILIAL", 2) | Out-Null

$d.Content.Find.Execute("Lugar: Indaial, Santa Catarina.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Lugar: Joinville, Santa Catarina.", 2) | Out-Null

# ---------------------------------------------------------------------------
# Step 3: insert the brand-new RADICALSOM.COM.BR / Premium block right
#         before the (now renamed) listing above.
# ---------------------------------------------------------------------------
$radicalNewText = New-BlockText "Storm 200" `
    "https://produto.mercadolivre.com.br/MLB-3452376800-carregador-jfa-storm-200a-144v-mais-completa-smart-cca-220v-_JM" `
    "Carregador Jfa Storm 200a 14,4v Mais Completa Smart Cca 220v" `
    "774.88" "845.87" "RADICALSOM.COM.BR" "Premium" `
    "Indaial, Santa Catarina."

$idx = Get-ParagraphIndexContaining("MLB-4156504566")
$modeloPara = $d.Paragraphs.Item($idx - 1)   # the "Modelo: Storm 200" line that starts this block
$modeloPara.Range.InsertBefore($radicalNewText)

# ---------------------------------------------------------------------------
# Step 4: delete the old "EVOX ELETRONICOS" block (the one that used to
#         directly follow "KGMICOMERCIAL"); it is being relocated to the
#         top of the document in step 5.
# ---------------------------------------------------------------------------
$idx = Get-ParagraphIndexContaining("MLB-3193872154")
$blockStartPara = $d.Paragraphs.Item($idx - 1)   # "Modelo: Storm 200" line
$blockEndPara = $d.Paragraphs.Item($idx - 1 + 9) # trailing blank separator paragraph
$delRange = $d.Range($blockStartPara.Range.Start, $blockEndPara.Range.End)
$delRange.Delete()

# ---------------------------------------------------------------------------
# Step 5: insert the (relocated) "EVOX ELETRONICOS" block at the very start
#         of the document.
# ---------------------------------------------------------------------------
$evoxText = New-BlockText "Storm 200" `
    "https://produto.mercadolivre.com.br/MLB-3193872154-fonte-carregador-automotivo-jfa-storm-200a-sci-bivolt-carro-_JM" `
    "Fonte Carregador Automotivo Jfa Storm 200a Sci Bivolt Carro" `
    "657.19" "805.59" "EVOX ELETRONICOS" "Clássico" `
    "Campo Limpo Paulista, São Paulo."

$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertBefore($evoxText)

Write-Host "Done. Final paragraph count:" $d.Paragraphs.Count
